
$d = $word.ActiveDocument

function Replace-Paragraph($startAnchor, $endAnchor, $newXml) {
    $rngStart = $d.Content.Duplicate
    $foundStart = $rngStart.Find.Execute($startAnchor)
    if (-not $foundStart) {
        throw "Start anchor not found: $startAnchor"
    }
    $rngEnd = $d.Content.Duplicate
    $rngEnd.Start = $rngStart.Start
    $foundEnd = $rngEnd.Find.Execute($endAnchor)
    if (-not $foundEnd) {
        throw "End anchor not found: $endAnchor"
    }
    $fullRng = $d.Range($rngStart.Start, $rngEnd.End)
    $fullRng.InsertXML($newXml)
}

# 1) First "Rest area - Lobby" paragraph -> "Coffee Break (Lobby); ..."
Replace-Paragraph 'Break (Rest area - Lobby' 'Additional rest area - Room 581)' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1B101F0B" w14:textId="7D48BD80" w:rsidR="624FD8E5" w:rsidRDefault="624FD8E5" w:rsidP="4D3864E2"><w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r w:rsidRPr="4D3864E2"><w:t>Break</w:t></w:r><w:r w:rsidR="00DB2722"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">Coffee </w:t></w:r><w:r><w:t>Break (Lobby</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>)</w:t></w:r><w:r><w:t>;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B34A31"><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Additional r</w:t></w:r><w:r w:rsidR="00DB2722"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">est area </w:t></w:r><w:r w:rsidR="00260516" w:rsidRPr="00260516"><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>- Room 581</w:t></w:r><w:r w:rsidR="00341066"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>)</w:t></w:r></w:p>'

# 2) Second "R" + "est area" paragraph -> "Coffee Break (Lobby); ..."
Replace-Paragraph 'Break (Rest area - Lobby' 'Additional rest area - Room 581)' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="78E350DF" w14:textId="31EF8A39" w:rsidR="4D3864E2" w:rsidRPr="000F2544" w:rsidRDefault="4D3864E2" w:rsidP="4D3864E2"><w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r w:rsidRPr="4D3864E2"><w:t>Break</w:t></w:r><w:r w:rsidR="00DB2722"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">Coffee </w:t></w:r><w:r><w:t>Break (Lobby</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>)</w:t></w:r><w:r><w:t>;</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B34A31"><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Additional r</w:t></w:r><w:r w:rsidR="00B34A31"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">est area - </w:t></w:r><w:r w:rsidR="00260516" w:rsidRPr="00260516"><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Room 581</w:t></w:r><w:r w:rsidR="000F2544"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>)</w:t></w:r></w:p>'

# 3) "Session 3A: Invited papers 1..." -> split "Invited paper" with proofErr tags
Replace-Paragraph 'Session 3A' '(Room 475)' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="33D88D1A" w14:textId="15DE6760" w:rsidR="00630C57" w:rsidRPr="009624DA" w:rsidRDefault="00630C57" w:rsidP="4D3864E2"><w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:t>Session 3A</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00243767"><w:t xml:space="preserve">:  </w:t></w:r><w:r w:rsidR="00243767" w:rsidRPr="00243767"><w:t>Invited</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> paper</w:t></w:r><w:r w:rsidR="008B37CA"><w:t>s 1</w:t></w:r><w:r w:rsidR="00243767" w:rsidRPr="00243767"><w:t>: AI/ML Modeling and Applications</w:t></w:r><w:r w:rsidR="00745661"><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> (Room 475)</w:t></w:r></w:p>'

# 4) "15:30 - 17:00" -> "15:30 - 16:40"
Replace-Paragraph '15:30 – 17:00' '15:30 – 17:00' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="106BA05C" w14:textId="4C8D2B9F" w:rsidR="005F641B" w:rsidRDefault="005F641B" w:rsidP="4D3864E2"><w:pPr><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr></w:pPr><w:r w:rsidRPr="4D3864E2"><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr><w:t>15:30 – 1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/></w:rPr><w:t>0</w:t></w:r></w:p>'

# 5) "DLS 2 + DT-PISCC  (Room 345)" -> split with proofErr tags
Replace-Paragraph 'DLS 2 + DT-PISCC' '(Room 345)' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="557B536B" w14:textId="77777777" w:rsidR="00EA6487" w:rsidRDefault="00EA6487" w:rsidP="00F673CA"><w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:t>DLS 2 + DT-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>PISCC</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">  (</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Room 345)</w:t></w:r></w:p>'

Write-Host "Done"
